$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 13:14 to make room for the two "Docentes responsaveis" rows
$ws.Rows("13:14").Insert()
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(14,1).Clear()

# --- Row 10/11: fill in real Objetivos/Objectives text (was filled with the ementa placeholder) ---
$ws.Cells.Item(10,2).Value = 'Fornecer meios para que o estudante adquira conhecimento e prática no uso de ferramentas computacionais modernas aplicadas a problemas físicos.'
$ws.Cells.Item(10,3).Value = 'Fornecer meios para que o estudante adquira conhecimento e prática no uso de ferramentas computacionais modernas aplicadas a problemas físicos.'
$ws.Cells.Item(11,2).Value = 'Provide means for the student to acquire knowledge and practice in the use of modern computational tools applied to physical problems.'
$ws.Cells.Item(11,3).Value = 'Provide means for the student to acquire knowledge and practice in the use of modern computational tools applied to physical problems.'

# --- Rows 13/14: Docentes responsaveis (two teachers) ---
$ws.Cells.Item(13,2).Value = '7290967 - Emerson Gonçalves de Melo'
$ws.Cells.Item(13,3).Value = '7290967 - Emerson Gonçalves de Melo'
$ws.Cells.Item(14,2).Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Cells.Item(14,3).Value = '1176388 - Luiz Tadeu Fernandes Eleno'
# Newly-inserted B13/B14 cells inherit column A bold/no-wrap style; restore the
# normal wrapped style used by every other column-B cell in this sheet.
$ws.Cells.Item(13,2).WrapText = $true
$ws.Cells.Item(13,2).Font.Bold = $false
$ws.Cells.Item(14,2).WrapText = $true
$ws.Cells.Item(14,2).Font.Bold = $false

# --- Row 15/16: Programa resumido / Short syllabus ---
$ws.Cells.Item(15,2).Value = 'Simulação numérica em sistemas determinísticos e estocásticos. Métodos de Monte Carlo. Caminhadas aleatórias. Fractais. Introdução à análise espectral por transformadas de Fourier. Revisão das soluções de equações diferenciais ordinárias e parciais. Solução numérica de equações diferenciais parciais.'
$ws.Cells.Item(15,3).Value = 'Simulação numérica em sistemas determinísticos e estocásticos. Métodos de Monte Carlo. Caminhadas aleatórias. Fractais. Introdução à análise espectral por transformadas de Fourier. Revisão das soluções de equações diferenciais ordinárias e parciais. Solução numérica de equações diferenciais parciais.'
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60

# --- Row 17/18: Programa / Syllabus ---
$ws.Cells.Item(17,2).Value = '• Simulação numérica em sistemas determinísticos. • Simulações numéricas em sistemas estocásticos. • Números pseudo-aleatórios. • O problema do caminho aleatório, difusão e percolação. • Fractais. • Introdução à análise espectral por transformadas de Fourier. • A Transformada rápida de Fourier e aplicações. • Identificação de frequências e modos normais. • Detecção e tratamento de sinais • Tratamento de imagens. • Métodos de solução numérica de equações diferenciais parciais. • Método das diferenças Finitas.  • Método dos Elementos Finitos.'
$ws.Cells.Item(17,3).Value = '• Simulação numérica em sistemas determinísticos. • Simulações numéricas em sistemas estocásticos. • Números pseudo-aleatórios. • O problema do caminho aleatório, difusão e percolação. • Fractais. • Introdução à análise espectral por transformadas de Fourier. • A Transformada rápida de Fourier e aplicações. • Identificação de frequências e modos normais. • Detecção e tratamento de sinais • Tratamento de imagens. • Métodos de solução numérica de equações diferenciais parciais. • Método das diferenças Finitas.  • Método dos Elementos Finitos.'
$ws.Rows.Item(17).RowHeight = 120

# --- Row 20: Metodo ---
$ws.Cells.Item(20,2).Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Cells.Item(20,3).Value = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: Criterio ---
$ws.Cells.Item(21,2).Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Cells.Item(21,3).Value = 'Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22: Norma de recuperacao ---
$ws.Cells.Item(22,2).Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Cells.Item(22,3).Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows.Item(22).RowHeight = 60

# --- Row 23: Bibliografia ---
$biblio = @'
GOULD, H.; TOBOCHNIK, T. An Introduction to Computer Simulation Methods. Addison-Wesley Publishing Company, Nova Iorque, 1987. 
SCHERER, C. Métodos Computacionais da Física, Editora Livraria da Física, São Paulo, 2005.
DEVRIES, P. L. A First Course in Computational Physics. John Wiley and Sons, New York, 1994.
PANG, H. An Introduction to Computational Physics. Cambridge University Press, Cambridge, 1997.
THIJSSEN, J. M. Computational Physics. Cambridge University Press, Cambridge, 1999. 
PRESS, W. H.; FLANNERY, B. P.; TEUKOLSKI, S. A.; VETERLING, W. T. Numerical Recipes. Cambridge University Press, 1986.
KOONIN, S. E. Computational Physics. Benjamin Cummings, 1986.
'@
$ws.Cells.Item(23,2).Value = $biblio
$ws.Cells.Item(23,3).Value = $biblio

